# Applies the figure_1_processing_workflow.pptx edit:
#  1. Refresh the cached "update automatically" date placeholder text
#     (6/29/21 -> 11/25/21) on the slide master and every slide layout.
#  2. Rename the "Variable extractions" bullet to "Descriptor extractions"
#     in the two code-listing boxes on slide 1.

$p = $ppt.ActivePresentation

# --- 1. Date placeholders -------------------------------------------------
$newDate = "11/25/21"

$master = $p.SlideMaster

for ($masterShapeIdx = 1; $masterShapeIdx -le $master.Shapes.Count; $masterShapeIdx++) {
    $masterShape = $master.Shapes.Item($masterShapeIdx)
    if ($masterShape.Name -like "Date Placeholder*") {
        $masterShape.TextFrame.TextRange.Text = $newDate
    }
}

for ($layoutIdx = 1; $layoutIdx -le $master.CustomLayouts.Count; $layoutIdx++) {
    $layout = $master.CustomLayouts.Item($layoutIdx)
    for ($layoutShapeIdx = 1; $layoutShapeIdx -le $layout.Shapes.Count; $layoutShapeIdx++) {
        $layoutShape = $layout.Shapes.Item($layoutShapeIdx)
        if ($layoutShape.Name -like "Date Placeholder*") {
            $layoutShape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2. "Variable extractions" -> "Descriptor extractions" ---------------
$oldLabel = "Variable extractions"
$newLabel = "Descriptor extractions"

$slide = $p.Slides.Item(1)

for ($slideShapeIdx = 1; $slideShapeIdx -le $slide.Shapes.Count; $slideShapeIdx++) {
    $slideShape = $slide.Shapes.Item($slideShapeIdx)
    if ($slideShape.HasTextFrame -eq -1) {
        $wholeText = $slideShape.TextFrame.TextRange.Text
        if ($wholeText -like "*$oldLabel*") {
            $shapeTextRange = $slideShape.TextFrame.TextRange
            $paraCount = $shapeTextRange.Paragraphs().Count
            for ($paraIdx = 1; $paraIdx -le $paraCount; $paraIdx++) {
                $paraRange = $shapeTextRange.Paragraphs($paraIdx, 1)
                if ($paraRange.Text.Trim() -eq $oldLabel) {
                    $runRange = $paraRange.Runs(1, 1)
                    $runRange.Text = $newLabel
                }
            }
        }
    }
}
